# "math jax working loop mode working"
#
# The three multiple-choice answer cells on row 7 of the "questions" sheet
# hold LaTeX-style math strings consumed by a MathJax-powered quiz
# renderer ($x^2 * x^3$, $x^5$, $x^2 + x^3$). To get MathJax's loop mode
# picking these up correctly they each need to be wrapped in an extra
# pair of curly braces, e.g. $x^5$ -> ${x^5}$.
#
# Update the three cells in place (A7, B7, C7) without disturbing
# anything else on the sheet (selection stays on C7, same as before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

$ws.Range("A7").Value2 = '${x^2 * x^3}$'
$ws.Range("B7").Value2 = '${x^5}$'
$ws.Range("C7").Value2 = '${x^2 + x^3}$'
